$wb = $excel.ActiveWorkbook

# Remember which sheet was active before we start, so the original selection
# state can be restored once the new sheet has been created.
$originalActive = $wb.ActiveSheet.Name

# Duplicate the "Summary" sheet to create a new "Summary_1" sheet right after it.
$summary = $wb.Worksheets.Item("Summary")
$summary.Copy($null, $summary)

$newSheet = $wb.Worksheets.Item($summary.Index + 1)
$newSheet.Name = "Summary_1"

# Apply percentage number formatting to the Percentage column (D2:D6).
$newSheet.Range("D2:D6").NumberFormat = "0.00%"

# Restore the originally active sheet/selection.
$wb.Worksheets.Item($originalActive).Activate()
